$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1054.8
$ws.Range("C2").Value = 995.6
$ws.Range("D2").Value = 1025.25
$ws.Range("E2").Value = 1028.1
$ws.Range("F2").Value = 187
$ws.Range("G2").Value = 1048.55

$ws.Range("B3").Value = 1037.2
$ws.Range("C3").Value = 1007.55
$ws.Range("D3").Value = 1022
$ws.Range("E3").Value = 1022.35
$ws.Range("F3").Value = 32
$ws.Range("G3").Value = 1033.7

$ws.Range("B4").Value = 47599
$ws.Range("C4").Value = 46961.3
$ws.Range("D4").Value = 47590
$ws.Range("E4").Value = 47499.8
$ws.Range("F4").Value = 27
$ws.Range("G4").Value = 47095

$ws.Range("B5").Value = 443.3
$ws.Range("C5").Value = 431.25
$ws.Range("D5").Value = 440.1
$ws.Range("E5").Value = 440.6
$ws.Range("F5").Value = 261
$ws.Range("G5").Value = 438.5

$ws.Range("B6").Value = 666
$ws.Range("C6").Value = 643.25
$ws.Range("D6").Value = 651.8
$ws.Range("E6").Value = 651.85
$ws.Range("F6").Value = 148
$ws.Range("G6").Value = 665

$ws.Range("B7").Value = 527.5
$ws.Range("C7").Value = 515.7
$ws.Range("D7").Value = 522.85
$ws.Range("E7").Value = 522.75
$ws.Range("F7").Value = 94
$ws.Range("G7").Value = 523.25

$ws.Range("B8").Value = 1017.95
$ws.Range("C8").Value = 999.5
$ws.Range("D8").Value = 1016.35
$ws.Range("E8").Value = 1015.9
$ws.Range("F8").Value = 222
$ws.Range("G8").Value = 1000.15

$ws.Range("B9").Value = 701.55
$ws.Range("C9").Value = 679.65
$ws.Range("D9").Value = 689.3
$ws.Range("E9").Value = 690.25
$ws.Range("F9").Value = 70
$ws.Range("G9").Value = 694.65

$ws.Range("B10").Value = 21107.45
$ws.Range("C10").Value = 20949
$ws.Range("D10").Value = 21090.2
$ws.Range("E10").Value = 21075
$ws.Range("F10").Value = 65
$ws.Range("G10").Value = 21071.1

$ws.Range("B11").Value = 2486.3
$ws.Range("C11").Value = 2450.7
$ws.Range("D11").Value = 2464
$ws.Range("E11").Value = 2465.05
$ws.Range("F11").Value = 63
$ws.Range("G11").Value = 2469.35

$ws.Range("B12").Value = 619.5
$ws.Range("C12").Value = 607.15
$ws.Range("D12").Value = 617.7
$ws.Range("E12").Value = 616.75
$ws.Range("F12").Value = 485
$ws.Range("G12").Value = 612.55

$ws.Range("B13").Value = 965.75
$ws.Range("C13").Value = 947.8
$ws.Range("D13").Value = 954.5
$ws.Range("E13").Value = 952.65
$ws.Range("F13").Value = 15
$ws.Range("G13").Value = 965.35

$ws.Range("B14").Value = 730.8
$ws.Range("C14").Value = 710.1
$ws.Range("D14").Value = 718.2
$ws.Range("E14").Value = 718.75
$ws.Range("F14").Value = 175
$ws.Range("G14").Value = 726.05

$ws.Range("B15").Value = 132.15
$ws.Range("C15").Value = 128.55
$ws.Range("D15").Value = 129.95
$ws.Range("E15").Value = 130
$ws.Range("F15").Value = 528
$ws.Range("G15").Value = 131.3

$ws.Range("B16").Value = 3662.05
$ws.Range("C16").Value = 3621.75
$ws.Range("D16").Value = 3643.9
$ws.Range("E16").Value = 3644.45
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 3635

$ws.Range("B17").Value = 3654.2
$ws.Range("C17").Value = 3600.1
$ws.Range("D17").Value = 3642
$ws.Range("E17").Value = 3645.2
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 3604.75
